$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("F44").Value = "GetParentTypeUseElements"
$ws.Range("F45").Value = "GetParentTypeUseElements"

$ws.Range("F8").Value = "ParseBlockElement"

$ws.Range("G8").Value = "Free-standing blocks done by ParseBlockElement. Most blocks parsed as components of e.g. IfStatements"

$ws.Range("F48").Value = "ParseMethodElement"
$ws.Range("F52").Value = "ParseMethodElement"
$ws.Range("F54").Value = "ParseMethodElement"

$ws.Range("G9").Select()
